# Omaha_Cal_Info_CP05MOAS-GL340_00001.xlsx
# "Added CP05MOAS-GL eng bar codes, corrected a couple reference designators"

$wb = $excel.ActiveWorkbook

$wsMoorings = $wb.Worksheets.Item("Moorings")
$wsCalInfo  = $wb.Worksheets.Item("Asset_Cal_Info")

# --- Correct the Sensor OOIBARCODE (column E) reference designators on the
#     Asset_Cal_Info sheet. Every row currently shares the bar code
#     "A00284" (shared string 35) - replace with the correct, per-group
#     ENG bar codes. New shared strings must be introduced in this order
#     so they land at the same indices as the target workbook:
#       38 OL000350, 39 A01759, 40 N00034, 41 N00033, 42 N00032, 43 N00031

# ADCPAM group (rows 2-5) -> N00031 (set last, see ordering note above)
# FLORTM group (rows 7-10) -> N00032
# CTDGVM group (row 12) -> N00033
# DOSTAM group (row 14) -> N00034
# PARADM group (row 16) -> A01759
# ENG group (row 18) -> OL000350

$wsCalInfo.Range("E18").Value = "OL000350"
$wsCalInfo.Range("E16").Value = "A01759"
$wsCalInfo.Range("E14").Value = "N00034"
$wsCalInfo.Range("E12").Value = "N00033"

$wsCalInfo.Range("E7").Value  = "N00032"
$wsCalInfo.Range("E8").Value  = "N00032"
$wsCalInfo.Range("E9").Value  = "N00032"
$wsCalInfo.Range("E10").Value = "N00032"

$wsCalInfo.Range("E2").Value = "N00031"
$wsCalInfo.Range("E3").Value = "N00031"
$wsCalInfo.Range("E4").Value = "N00031"
$wsCalInfo.Range("E5").Value = "N00031"

# --- Restore the selection left behind on each sheet by the editor ---
$wsMoorings.Activate()
$wsMoorings.Range("K8:K9").Select()

$wsCalInfo.Activate()
$wsCalInfo.Range("H15").Select()
